$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New dictionary entries: Spanish names for the moon-phase rows (column D) ---
$ws.Range("D93").Value = "menguante gibosa"
$ws.Range("D94").Value = "cuarto menguante"
$ws.Range("D95").Value = "menguante"
$ws.Range("D96").Value = "nueva"
$ws.Range("D97").Value = "creciente"
$ws.Range("D98").Value = "cuarto creciente"
$ws.Range("D99").Value = "menguant gibosa"

# --- View state: zoom + active selection now on column D ---
$excel.ActiveWindow.Zoom = 150
$ws.Range("D99").Select()

# --- Page setup: explicit portrait orientation ---
$ws.PageSetup.Orientation = 1
